# Delete row 2 of the active sheet ("5"); this shifts rows 3..9 up to become rows 2..8,
# updates the merged cell ranges (A5:A6/M5:M6 -> A4:A5/M4:M5, A7:A8/M7:M8 -> A6:A7/M6:M7),
# the SUM formulas (SUM(J5:J6)->SUM(J4:J5), SUM(J7:J8)->SUM(J6:J7)) and the hyperlinks
# (E9..I9 -> E8..I8) automatically, matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(2).Delete()

# Update the "strategy 2 comment" text in the K2 cell of worksheet "8" to a new,
# distinct comment ("strategy 3 comment"), which becomes a new shared string.
$ws8 = $wb.Worksheets.Item("8")
$ws8.Range("K2").Value = "strategy 3 comment"
